# Applies the cryptos.xlsx price/volume refresh described by the commit:
# "Updated cryptos list on Mon Oct 23 09:14:00 UTC 2023 with GitHub Actions"
#
# Coin/Link/Price/Volume(1h) columns are plain text in the source sheet (no
# numeric typing), so every write goes through Set-TextValue, which forces the
# cell to Text format before assigning, then restores the default "Normal"
# style so number-looking strings (e.g. "220.35") land as text instead of
# being auto-coerced into a Double by the COM Value setter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '30.597.03'
Set-TextValue 'E2' '  +2.23%  '
# Row 3
Set-TextValue 'D3' '1.678.74'
Set-TextValue 'E3' '  +2.82%  '
# Row 4
Set-TextValue 'E4' '  -0.04%  '
# Row 5
Set-TextValue 'D5' '220.35'
Set-TextValue 'E5' '  +2.75%  '
# Row 6
Set-TextValue 'D6' '0.533'
Set-TextValue 'E6' '  +3.07%  '
# Row 7
Set-TextValue 'E7' '  -0.01%  '
# Row 8
Set-TextValue 'D8' '30.12'
Set-TextValue 'E8' '  +5.08%  '
# Row 9
Set-TextValue 'D9' '0.264'
Set-TextValue 'E9' '  +2.57%  '
# Row 10
Set-TextValue 'D10' '0.0640'
Set-TextValue 'E10' '  +5.21%  '
# Row 11
Set-TextValue 'E11' '  -0.67%  '
# Row 12
Set-TextValue 'D12' '1.921.01'
Set-TextValue 'E12' '  +2.93%  '
# Row 13
Set-TextValue 'D13' '10.41'
Set-TextValue 'E13' '  +12.03%  '
# Row 14
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.680.48'
Set-TextValue 'E14' '  +2.89%  '
# Row 15
Set-TextValue 'B15' 'Polygon'
Set-TextValue 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.616'
Set-TextValue 'E15' '  +9.18%  '
# Row 16
Set-TextValue 'E16' '  +3.44%  '
# Row 17
Set-TextValue 'D17' '30.607.12'
Set-TextValue 'E17' '  +2.33%  '
# Row 18
Set-TextValue 'D18' '66.49'
Set-TextValue 'E18' '  +3.72%  '
# Row 19
Set-TextValue 'D19' '245.67'
Set-TextValue 'E19' '  +1.17%  '
# Row 20
Set-TextValue 'D20' '0.0₃0727'
Set-TextValue 'E20' '  +3.73%  '
# Row 21
Set-TextValue 'D21' '0.998'
Set-TextValue 'E21' '  -0.25%  '
# Row 22
Set-TextValue 'D22' '4.28'
Set-TextValue 'E22' '  +3.73%  '
# Row 23
Set-TextValue 'D23' '10.08'
Set-TextValue 'E23' '  +2.70%  '
# Row 24
Set-TextValue 'D24' '2.16'
Set-TextValue 'E24' '  +1.07%  '
# Row 25
Set-TextValue 'D25' '158.33'
Set-TextValue 'E25' '  +0.27%  '
# Row 26
Set-TextValue 'D26' '15.95'
Set-TextValue 'E26' '  +2.80%  '
# Row 27
Set-TextValue 'D27' '0.113'
Set-TextValue 'E27' '  +2.68%  '
# Row 28
Set-TextValue 'D28' '6.71'
Set-TextValue 'E28' '  +1.97%  '
# Row 29
Set-TextValue 'E29' '  -0.04%  '
# Row 30
Set-TextValue 'D30' '0.0497'
Set-TextValue 'E30' '  +2.46%  '
# Row 31
Set-TextValue 'E31' '  +3.08%  '
# Row 32
Set-TextValue 'D32' '3.49'
Set-TextValue 'E32' '  +3.37%  '
# Row 33
Set-TextValue 'D33' '1.511.74'
Set-TextValue 'E33' '  +6.23%  '
# Row 34
Set-TextValue 'D34' '3.30'
Set-TextValue 'E34' '  +4.45%  '
# Row 35
Set-TextValue 'E35' '  +7.79%  '
# Row 36
Set-TextValue 'E36' '  -0.25%  '
# Row 37
Set-TextValue 'D37' '83.88'
Set-TextValue 'E37' '  +11.60%  '
# Row 38
Set-TextValue 'B38' 'ImmutableX'
Set-TextValue 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D38' '0.600'
Set-TextValue 'E38' '  +8.76%  '
# Row 39
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.0179'
Set-TextValue 'E39' '  +5.57%  '
# Row 40
Set-TextValue 'D40' '2.71'
Set-TextValue 'E40' '  -3.41%  '
# Row 41
Set-TextValue 'E41' '  +0.38%  '
# Row 42
Set-TextValue 'D42' '0.842'
Set-TextValue 'E42' '  +1.84%  '
# Row 43
Set-TextValue 'E43' '  +0.18%  '
# Row 44
Set-TextValue 'E44' '  +1.97%  '
# Row 45
Set-TextValue 'E45' '  +0.77%  '
# Row 46
Set-TextValue 'E46' '  -0.10%  '
# Row 47
Set-TextValue 'B47' 'FraxShare'
Set-TextValue 'C47' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D47' '5.62'
Set-TextValue 'E47' '  +5.26%  '
# Row 48
Set-TextValue 'B48' 'BitcoinSV'
Set-TextValue 'C48' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D48' '51.74'
Set-TextValue 'E48' '  -2.36%  '
# Row 49
Set-TextValue 'D49' '1.815.23'
# Row 50
Set-TextValue 'D50' '95.17'
Set-TextValue 'E50' '  +6.74%  '
# Row 51
Set-TextValue 'D51' '0.0₆0112'
Set-TextValue 'E51' '  +0.75%  '
